$d = $word.ActiveDocument

# 1. Fix wording: "brainstorming" -> "brainstormed" and add "Sir Paolo" after "adviser".
$d.Content.Find.Execute(
    "I’ve brainstorming about what App I will be developing. I came up with three app proposals and end up choosing the Music app based on the advice of my adviser. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I’ve brainstormed about what App I will be developing. I came up with three app proposals and end up choosing the Music app based on the advice of my adviser Sir Paolo. ",
    2
) | Out-Null

# 2. Move the "_GoBack" bookmark from the near-empty paragraph at the end of the
#    document to the start of the "Music App" paragraph (right before its first run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Paragraphs(4).Range
$target.Collapse(1)
$d.Bookmarks.Add("_GoBack", $target) | Out-Null
